$d = $word.ActiveDocument

$d.Content.Find.Execute("335×9=", $true, $false, $false, $false, $false, $true, 1, $false, "110×7=", 2) | Out-Null
$d.Content.Find.Execute("122×9=", $true, $false, $false, $false, $false, $true, 1, $false, "324×2=", 2) | Out-Null
$d.Content.Find.Execute("339×2=", $true, $false, $false, $false, $false, $true, 1, $false, "187×8=", 2) | Out-Null
$d.Content.Find.Execute("953×4=", $true, $false, $false, $false, $false, $true, 1, $false, "159×7=", 2) | Out-Null
$d.Content.Find.Execute("889×6=", $true, $false, $false, $false, $false, $true, 1, $false, "903×2=", 2) | Out-Null
$d.Content.Find.Execute("865×7=", $true, $false, $false, $false, $false, $true, 1, $false, "578×5=", 2) | Out-Null
$d.Content.Find.Execute("470×5=", $true, $false, $false, $false, $false, $true, 1, $false, "828×6=", 2) | Out-Null
$d.Content.Find.Execute("736×5=", $true, $false, $false, $false, $false, $true, 1, $false, "932×9=", 2) | Out-Null
$d.Content.Find.Execute("649×6=", $true, $false, $false, $false, $false, $true, 1, $false, "502×8=", 2) | Out-Null
$d.Content.Find.Execute("128×2=", $true, $false, $false, $false, $false, $true, 1, $false, "231×3=", 2) | Out-Null
$d.Content.Find.Execute("112×7=", $true, $false, $false, $false, $false, $true, 1, $false, "199×7=", 2) | Out-Null
$d.Content.Find.Execute("214×8=", $true, $false, $false, $false, $false, $true, 1, $false, "852×4=", 2) | Out-Null
$d.Content.Find.Execute("239×9=", $true, $false, $false, $false, $false, $true, 1, $false, "840×6=", 2) | Out-Null
$d.Content.Find.Execute("984×8=", $true, $false, $false, $false, $false, $true, 1, $false, "414×2=", 2) | Out-Null
$d.Content.Find.Execute("746×9=", $true, $false, $false, $false, $false, $true, 1, $false, "625×3=", 2) | Out-Null
$d.Content.Find.Execute("358×5=", $true, $false, $false, $false, $false, $true, 1, $false, "596×3=", 2) | Out-Null
$d.Content.Find.Execute("821×6=", $true, $false, $false, $false, $false, $true, 1, $false, "589×9=", 2) | Out-Null
$d.Content.Find.Execute("862×4=", $true, $false, $false, $false, $false, $true, 1, $false, "745×3=", 2) | Out-Null
$d.Content.Find.Execute("309×4=", $true, $false, $false, $false, $false, $true, 1, $false, "715×6=", 2) | Out-Null
$d.Content.Find.Execute("490×8=", $true, $false, $false, $false, $false, $true, 1, $false, "400×2=", 2) | Out-Null
$d.Content.Find.Execute("390×8=", $true, $false, $false, $false, $false, $true, 1, $false, "320×2=", 2) | Out-Null
$d.Content.Find.Execute("140×3=", $true, $false, $false, $false, $false, $true, 1, $false, "827×6=", 2) | Out-Null
$d.Content.Find.Execute("539×6=", $true, $false, $false, $false, $false, $true, 1, $false, "303×7=", 2) | Out-Null
$d.Content.Find.Execute("984×4=", $true, $false, $false, $false, $false, $true, 1, $false, "873×6=", 2) | Out-Null
$d.Content.Find.Execute("338×3=", $true, $false, $false, $false, $false, $true, 1, $false, "269×8=", 2) | Out-Null
